$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 12807.75
$ws.Range("J32").Value = 12807.75
$ws.Range("L32").Value = 12807.75
$ws.Range("N32").Value = -13459.75

$ws.Range("H40").Value = 1256628.2
$ws.Range("I40").Value = 1674918.5
$ws.Range("K40").Value = 1674918.5
$ws.Range("M40").Value = -1674743.5

$ws.Range("H53").Value = 1165.9412
$ws.Range("I53").Value = 1108.8182
$ws.Range("J53").Value = 1270.6666
$ws.Range("K53").Value = 1108.8182
$ws.Range("L53").Value = 1270.6666
$ws.Range("M53").Value = -471.8181999999999
$ws.Range("N53").Value = -2544.6666

$ws.Range("H62").Value = 4161.1665
$ws.Range("I62").Value = 3993.4
$ws.Range("K62").Value = 3993.4
$ws.Range("M62").Value = -3369.4

$ws.Range("H65").Value = 4161.1665
$ws.Range("I65").Value = 3993.4
$ws.Range("K65").Value = 19967
$ws.Range("M65").Value = -16847

$ws.Range("H74").Value = 5538.769
$ws.Range("I74").Value = 3002.5
$ws.Range("J74").Value = 5999.909
$ws.Range("K74").Value = 3002.5
$ws.Range("L74").Value = 5999.909
$ws.Range("M74").Value = -2066.5
$ws.Range("N74").Value = -7871.909

$ws.Range("H77").Value = 5538.769
$ws.Range("I77").Value = 3002.5
$ws.Range("J77").Value = 5999.909
$ws.Range("K77").Value = 15012.5
$ws.Range("L77").Value = 29999.545
$ws.Range("M77").Value = -10332.5
$ws.Range("N77").Value = -39359.545

$ws.Range("H94").Value = 927
$ws.Range("I94").Value = 927
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 927
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -476
$ws.Range("N94").ClearContents()

$ws.Range("H116").Value = 8595.556
$ws.Range("I116").Value = 8519.4
$ws.Range("J116").Value = 8690.75
$ws.Range("K116").Value = 8519.4
$ws.Range("L116").Value = 8690.75
$ws.Range("M116").Value = -5077.4
$ws.Range("N116").Value = -15574.75

$ws.Range("H125").Value = 2836.125
$ws.Range("J125").Value = 3114.8333
$ws.Range("L125").Value = 28033.4997
$ws.Range("N125").Value = -32953.4997

$ws.Range("H127").Value = 1397.8462
$ws.Range("I127").Value = 894.4545000000001
$ws.Range("J127").Value = 4166.5
$ws.Range("K127").Value = 2683.3635
$ws.Range("L127").Value = 12499.5
$ws.Range("M127").Value = 2276.6365
$ws.Range("N127").Value = -22419.5

$ws.Range("H131").Value = 5328.143
$ws.Range("I131").Value = 3587.6667
$ws.Range("J131").Value = 6633.5
$ws.Range("K131").Value = 10763.0001
$ws.Range("L131").Value = 19900.5
$ws.Range("M131").Value = -5723.000100000001
$ws.Range("N131").Value = -29980.5

$ws.Range("H138").Value = 6804562.5
$ws.Range("I138").Value = 679.26086
$ws.Range("J138").Value = 12823383
$ws.Range("K138").Value = 2037.78258
$ws.Range("L138").Value = 38470149
$ws.Range("M138").Value = 3102.21742
$ws.Range("N138").Value = -38480429

$ws.Range("H141").Value = 1704.36
$ws.Range("I141").Value = 1582.5652
$ws.Range("K141").Value = 4747.6956
$ws.Range("M141").Value = 432.3044

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 6775.3335
$ws.Range("J24").Value = 6775.3335
$ws.Range("L24").Value = 6775.3335
$ws.Range("N24").Value = -7523.3335

$ws.Range("H45").Value = 6054.778
$ws.Range("I45").Value = 5096
$ws.Range("K45").Value = 5096
$ws.Range("M45").Value = -4719

$ws.Range("H63").Value = 6582.9165
$ws.Range("I63").Value = 3330
$ws.Range("J63").Value = 7667.222
$ws.Range("K63").Value = 3330
$ws.Range("L63").Value = 7667.222
$ws.Range("M63").Value = -2644
$ws.Range("N63").Value = -9039.222

$ws.Range("H66").Value = 6582.9165
$ws.Range("I66").Value = 3330
$ws.Range("J66").Value = 7667.222
$ws.Range("K66").Value = 16650
$ws.Range("L66").Value = 38336.11
$ws.Range("M66").Value = -13218
$ws.Range("N66").Value = -45200.11

$ws.Range("H100").Value = 6775.3335
$ws.Range("J100").Value = 6775.3335
$ws.Range("L100").Value = 6775.3335
$ws.Range("N100").Value = -8939.333500000001

$ws.Range("H132").Value = 27818570
$ws.Range("I132").Value = 2052.742
$ws.Range("J132").Value = 200280980
$ws.Range("K132").Value = 6158.226000000001
$ws.Range("L132").Value = 600842940
$ws.Range("M132").Value = -3628.226000000001
$ws.Range("N132").Value = -600848000

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 869.75
$ws.Range("I80").Value = 462.33334
$ws.Range("K80").Value = 462.33334
$ws.Range("M80").Value = 535.66666

$ws.Range("H83").Value = 869.75
$ws.Range("I83").Value = 462.33334
$ws.Range("K83").Value = 2311.6667
$ws.Range("M83").Value = 2680.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 210.44
$ws.Range("I7").Value = 52.53846
$ws.Range("J7").Value = 381.5
$ws.Range("K7").Value = 52.53846
$ws.Range("L7").Value = 381.5
$ws.Range("M7").Value = 60.46154
$ws.Range("N7").Value = -607.5

$ws.Range("H88").Value = 28724.75
$ws.Range("J88").Value = 28724.75
$ws.Range("L88").Value = 28724.75
$ws.Range("N88").Value = -29536.75

$ws.Range("H91").Value = 28724.75
$ws.Range("J91").Value = 28724.75
$ws.Range("L91").Value = 28724.75
$ws.Range("N91").Value = -31532.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 283.6
$ws.Range("I8").Value = 283.6
$ws.Range("K8").Value = 850.8000000000001
$ws.Range("M8").Value = -711.8000000000001

$ws.Range("H22").Value = 414.4375
$ws.Range("I22").Value = 148.27272
$ws.Range("K22").Value = 444.81816
$ws.Range("M22").Value = -275.81816

$ws.Range("H26").Value = 380.875
$ws.Range("I26").Value = 70
$ws.Range("J26").Value = 567.4
$ws.Range("K26").Value = 210
$ws.Range("L26").Value = 1702.2
$ws.Range("M26").Value = 78
$ws.Range("N26").Value = -2278.2

$ws.Range("H27").Value = 414.4375
$ws.Range("I27").Value = 148.27272
$ws.Range("K27").Value = 444.81816
$ws.Range("M27").Value = -342.81816

$ws.Range("H55").Value = 6667068.5
$ws.Range("I55").Value = 442.7143
$ws.Range("J55").Value = 22222528
$ws.Range("K55").Value = 1328.1429
$ws.Range("L55").Value = 66667584
$ws.Range("M55").Value = -1151.1429
$ws.Range("N55").Value = -66667938

$ws.Range("H88").Value = 20000
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 20000
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 60000
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -60856

$ws.Range("H91").Value = 20000
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 20000
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 60000
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -62964

$ws.Range("H113").Value = 2364.7727
$ws.Range("I113").Value = 964.2857
$ws.Range("K113").Value = 2892.8571
$ws.Range("M113").Value = -722.8571000000002

$ws.Range("H129").Value = 3415.6177
$ws.Range("I129").Value = 3579.8
$ws.Range("J129").Value = 3347.2083
$ws.Range("K129").Value = 10739.4
$ws.Range("L129").Value = 10041.6249
$ws.Range("M129").Value = -5739.400000000001
$ws.Range("N129").Value = -20041.6249

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 12820674
$ws.Range("I2").Value = 15873078
$ws.Range("J2").Value = 577.4
$ws.Range("K2").Value = 15873078
$ws.Range("L2").Value = 577.4
$ws.Range("M2").Value = -15872965
$ws.Range("N2").Value = -803.4

$ws.Range("H14").Value = 5474512.5
$ws.Range("I14").Value = 5474512.5
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 5474512.5
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -5474344.5
$ws.Range("N14").ClearContents()

$ws.Range("H70").Value = 3518.577
$ws.Range("I70").Value = 3566.1428
$ws.Range("K70").Value = 3566.1428
$ws.Range("M70").Value = -3296.1428

$ws.Range("H73").Value = 3518.577
$ws.Range("I73").Value = 3566.1428
$ws.Range("K73").Value = 3566.1428
$ws.Range("M73").Value = -2630.1428

$ws.Range("H93").Value = 46242.5
$ws.Range("J93").Value = 46242.5
$ws.Range("L93").Value = 46242.5
$ws.Range("N93").Value = -49986.5

$ws.Range("H132").Value = 1680.6154
$ws.Range("I132").Value = 1518.7778
$ws.Range("J132").Value = 2044.75
$ws.Range("K132").Value = 4556.3334
$ws.Range("L132").Value = 6134.25
$ws.Range("M132").Value = -2026.3334
$ws.Range("N132").Value = -11194.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 769.4761999999999
$ws.Range("I16").Value = 445.15384
$ws.Range("K16").Value = 445.15384
$ws.Range("M16").Value = -275.15384

$ws.Range("H50").Value = 25000
$ws.Range("J50").Value = 25000
$ws.Range("L50").Value = 25000
$ws.Range("N50").Value = -26274

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6621.4287
$ws.Range("I62").Value = 4925.5
$ws.Range("K62").Value = 4925.5
$ws.Range("M62").Value = -4301.5

$ws.Range("H65").Value = 6621.4287
$ws.Range("I65").Value = 4925.5
$ws.Range("K65").Value = 24627.5
$ws.Range("M65").Value = -21507.5

$ws.Range("H96").Value = 1799.75
$ws.Range("I96").Value = 1799.75
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 1799.75
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -426.75
$ws.Range("N96").ClearContents()

$ws.Range("H132").Value = 5458.091
$ws.Range("I132").Value = 5874
$ws.Range("J132").Value = 1299
$ws.Range("K132").Value = 17622
$ws.Range("L132").Value = 3897
$ws.Range("M132").Value = -15092
$ws.Range("N132").Value = -8957
